$d = $word.ActiveDocument
$wns = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ---------------------------------------------------------------------
# 1) "Headbob Mechanic" -> split into "Headbob" (wrapped in a
#    spellStart/spellEnd proofErr pair) + " Mechanic", leaving the
#    trailing "." run that follows it untouched (but reconstructed,
#    since InsertXML always lands its new content at the end of the
#    paragraph it targets).
# ---------------------------------------------------------------------
$fr = $d.Content
$fr.Find.Execute("Headbob Mechanic") | Out-Null
$s = $fr.Start
$para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $s -and $p.Range.End -gt $s) { $para = $p; break }
}
$rng = $d.Range($s, $para.Range.End)
$xml = '<w:p xmlns:w="' + $wns + '">' +
       '<w:proofErr w:type="spellStart"/>' +
       '<w:r><w:t>Headbob</w:t></w:r>' +
       '<w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> Mechanic</w:t></w:r>' +
       '<w:r><w:t>.</w:t></w:r>' +
       '</w:p>'
$rng.InsertXML($xml)

# ---------------------------------------------------------------------
# 2) "New GameMode" / "New GameInstance" / "New GameState" /
#    "New PlayerState" -> each single-run paragraph is split into a
#    "New " run plus a spellStart/spellEnd-wrapped run for the
#    CamelCase word. Each of these paragraphs holds exactly one run
#    that fills the whole paragraph, so a full-paragraph replace is
#    safe and keeps things in order.
# ---------------------------------------------------------------------
$words = @("GameMode", "GameInstance", "GameState", "PlayerState")
foreach ($w in $words) {
    $needle = "New " + $w
    $fr2 = $d.Content
    $fr2.Find.Execute($needle) | Out-Null
    $s2 = $fr2.Start
    $para2 = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p2 = $d.Paragraphs.Item($i)
        if ($p2.Range.Start -le $s2 -and $p2.Range.End -gt $s2) { $para2 = $p2; break }
    }
    $rng2 = $para2.Range.Duplicate
    $xml2 = '<w:p xmlns:w="' + $wns + '">' +
            '<w:r><w:t xml:space="preserve">New </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:r><w:t>' + $w + '</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '</w:p>'
    $rng2.InsertXML($xml2)
}

# ---------------------------------------------------------------------
# 3) Remove the duplicate "Bug causing drop down menu to appear
#    incorrectly" list item (the "display incorrectly" one right
#    before it is kept).
# ---------------------------------------------------------------------
$fr3 = $d.Content
$fr3.Find.Execute("Bug causing drop down menu to appear incorrectly") | Out-Null
$s3 = $fr3.Start
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p3 = $d.Paragraphs.Item($i)
    if ($p3.Range.Start -le $s3 -and $p3.Range.End -gt $s3) {
        $p3.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 4) Footer: wrap "HOWLONG" / "(" in a gramStart/gramEnd proofErr pair
#    (splitting the "(v" run into "(" and "v"), and refresh the DATE
#    field's cached display text. The footer paragraph contains a
#    DATE field (fldChar begin/separate/end) that this engine mishandles
#    under partial-range InsertXML, so the whole paragraph is rebuilt
#    from its start in one shot to keep the field intact and in order.
# ---------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$fpara = $footer.Range.Paragraphs.Item(1)
$frng = $fpara.Range.Duplicate
$fxml = '<w:p xmlns:w="' + $wns + '">' +
        '<w:r><w:t xml:space="preserve">Project </w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:t>HOWLONG</w:t></w:r>' +
        '<w:r><w:t>(</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t>v</w:t></w:r>' +
        '<w:r><w:t>3</w:t></w:r>' +
        '<w:r><w:t>.0)</w:t></w:r>' +
        '<w:r><w:t>20</w:t></w:r>' +
        '<w:r><w:t>2</w:t></w:r>' +
        '<w:r><w:t>2</w:t></w:r>' +
        '<w:r><w:t>-01</w:t></w:r>' +
        '<w:r><w:tab/></w:r>' +
        '<w:r><w:tab/><w:t xml:space="preserve">Last Updated </w:t></w:r>' +
        '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
        '<w:r><w:instrText xml:space="preserve"> DATE \@ &quot;dd/MM/yyyy HH:mm&quot; </w:instrText></w:r>' +
        '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
        '<w:r><w:rPr><w:noProof/></w:rPr><w:t>07/07/2025 15:24</w:t></w:r>' +
        '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
        '</w:p>'
$frng.InsertXML($fxml)

Write-Output "done"
